$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Apoe"
$ws.Range("C2").Value = "Vldlr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 37.05583833333333
$ws.Range("H2").Value = 111.167515
$ws.Range("I2").Value = 0.008431126118266585
$ws.Range("J2").Value = 0.008431126118266585
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3081963333333333
$ws.Range("N2").Value = 0.924589
$ws.Range("O2").Value = 0.09210955608663024
$ws.Range("P2").Value = 0.09210955608663024
$ws.Range("Q2").Value = 11.42047350292611
$ws.Range("R2").Value = 102.784261526335
$ws.Range("S2").Value = 0.0007765872840639291
$ws.Range("T2").Value = 0.0007765872840639291
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Apoe"
$ws.Range("C3").Value = "Vldlr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 37.05583833333333
$ws.Range("H3").Value = 111.167515
$ws.Range("I3").Value = 0.008431126118266585
$ws.Range("J3").Value = 0.008431126118266585
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.207039333333333
$ws.Range("N3").Value = 6.621118
$ws.Range("O3").Value = 0.65961009678592
$ws.Range("P3").Value = 0.6596100967859201
$ws.Range("Q3").Value = 81.78369273130778
$ws.Range("R3").Value = 736.05323458177
$ws.Range("S3").Value = 0.00556125591488412
$ws.Range("T3").Value = 0.005561255914884121
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Apoe"
$ws.Range("C4").Value = "Vldlr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 37.05583833333333
$ws.Range("H4").Value = 111.167515
$ws.Range("I4").Value = 0.008431126118266585
$ws.Range("J4").Value = 0.008431126118266585
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.83074
$ws.Range("N4").Value = 2.49222
$ws.Range("O4").Value = 0.2482803471274497
$ws.Range("P4").Value = 0.2482803471274497
$ws.Range("Q4").Value = 30.78376713703334
$ws.Range("R4").Value = 277.0539042333
$ws.Range("S4").Value = 0.002093282919318535
$ws.Range("T4").Value = 0.002093282919318535
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Apoe"
$ws.Range("C5").Value = "Vldlr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 50.89916233333333
$ws.Range("H5").Value = 152.697487
$ws.Range("I5").Value = 0.01158082710438721
$ws.Range("J5").Value = 0.01158082710438721
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3081963333333333
$ws.Range("N5").Value = 0.924589
$ws.Range("O5").Value = 0.09210955608663024
$ws.Range("P5").Value = 0.09210955608663024
$ws.Range("Q5").Value = 15.68693520087144
$ws.Range("R5").Value = 141.182416807843
$ws.Range("S5").Value = 0.001066704843701122
$ws.Range("T5").Value = 0.001066704843701122
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Apoe"
$ws.Range("C6").Value = "Vldlr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 50.89916233333333
$ws.Range("H6").Value = 152.697487
$ws.Range("I6").Value = 0.01158082710438721
$ws.Range("J6").Value = 0.01158082710438721
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.207039333333333
$ws.Range("N6").Value = 6.621118
$ws.Range("O6").Value = 0.65961009678592
$ws.Range("P6").Value = 0.6596100967859201
$ws.Range("Q6").Value = 112.3364533033851
$ws.Range("R6").Value = 1011.028079730466
$ws.Range("S6").Value = 0.007638830487185856
$ws.Range("T6").Value = 0.007638830487185857
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Apoe"
$ws.Range("C7").Value = "Vldlr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 50.89916233333333
$ws.Range("H7").Value = 152.697487
$ws.Range("I7").Value = 0.01158082710438721
$ws.Range("J7").Value = 0.01158082710438721
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.83074
$ws.Range("N7").Value = 2.49222
$ws.Range("O7").Value = 0.2482803471274497
$ws.Range("P7").Value = 0.2482803471274497
$ws.Range("Q7").Value = 42.28397011679333
$ws.Range("R7").Value = 380.55573105114
$ws.Range("S7").Value = 0.002875291773500236
$ws.Range("T7").Value = 0.002875291773500236
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Apoe"
$ws.Range("C8").Value = "Vldlr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2007.446289
$ws.Range("H8").Value = 6022.338867
$ws.Range("I8").Value = 0.4567440273772037
$ws.Range("J8").Value = 0.4567440273772037
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.3081963333333333
$ws.Range("N8").Value = 0.924589
$ws.Range("O8").Value = 0.09210955608663024
$ws.Range("P8").Value = 0.09210955608663024
$ws.Range("Q8").Value = 618.6875856334071
$ws.Range("R8").Value = 5568.188270700663
$ws.Range("S8").Value = 0.04207048960693392
$ws.Range("T8").Value = 0.04207048960693392
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Apoe"
$ws.Range("C9").Value = "Vldlr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2007.446289
$ws.Range("H9").Value = 6022.338867
$ws.Range("I9").Value = 0.4567440273772037
$ws.Range("J9").Value = 0.4567440273772037
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.207039333333333
$ws.Range("N9").Value = 6.621118
$ws.Range("O9").Value = 0.65961009678592
$ws.Range("P9").Value = 0.6596100967859201
$ws.Range("Q9").Value = 4430.512919377034
$ws.Range("R9").Value = 39874.61627439331
$ws.Range("S9").Value = 0.3012729721046682
$ws.Range("T9").Value = 0.3012729721046682
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Apoe"
$ws.Range("C10").Value = "Vldlr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2007.446289
$ws.Range("H10").Value = 6022.338867
$ws.Range("I10").Value = 0.4567440273772037
$ws.Range("J10").Value = 0.4567440273772037
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.83074
$ws.Range("N10").Value = 2.49222
$ws.Range("O10").Value = 0.2482803471274497
$ws.Range("P10").Value = 0.2482803471274497
$ws.Range("Q10").Value = 1667.66593012386
$ws.Range("R10").Value = 15008.99337111474
$ws.Range("S10").Value = 0.1134005656656015
$ws.Range("T10").Value = 0.1134005656656015
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Apoe"
$ws.Range("C11").Value = "Vldlr"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2293.273345666667
$ws.Range("H11").Value = 6879.820037
$ws.Range("I11").Value = 0.5217768014597114
$ws.Range("J11").Value = 0.5217768014597114
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.3081963333333333
$ws.Range("N11").Value = 0.924589
$ws.Range("O11").Value = 0.09210955608663024
$ws.Range("P11").Value = 0.09210955608663024
$ws.Range("Q11").Value = 706.7784364655326
$ws.Range("R11").Value = 6361.005928189793
$ws.Range("S11").Value = 0.04806062955875582
$ws.Range("T11").Value = 0.04806062955875582
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Apoe"
$ws.Range("C12").Value = "Vldlr"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2293.273345666667
$ws.Range("H12").Value = 6879.820037
$ws.Range("I12").Value = 0.5217768014597114
$ws.Range("J12").Value = 0.5217768014597114
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.207039333333333
$ws.Range("N12").Value = 6.621118
$ws.Range("O12").Value = 0.65961009678592
$ws.Range("P12").Value = 0.6596100967859201
$ws.Range("Q12").Value = 5061.344475971264
$ws.Range("R12").Value = 45552.10028374137
$ws.Range("S12").Value = 0.344169246511488
$ws.Range("T12").Value = 0.3441692465114881
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Apoe"
$ws.Range("C13").Value = "Vldlr"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2293.273345666667
$ws.Range("H13").Value = 6879.820037
$ws.Range("I13").Value = 0.5217768014597114
$ws.Range("J13").Value = 0.5217768014597114
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.83074
$ws.Range("N13").Value = 2.49222
$ws.Range("O13").Value = 0.2482803471274497
$ws.Range("P13").Value = 0.2482803471274497
$ws.Range("Q13").Value = 1905.113899179127
$ws.Range("R13").Value = 17146.02509261214
$ws.Range("S13").Value = 0.1295469253894675
$ws.Range("T13").Value = 0.1295469253894676
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Apoe"
$ws.Range("C14").Value = "Vldlr"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 6.448603666666666
$ws.Range("H14").Value = 19.345811
$ws.Range("I14").Value = 0.00146721794043115
$ws.Range("J14").Value = 0.00146721794043115
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.3081963333333333
$ws.Range("N14").Value = 0.924589
$ws.Range("O14").Value = 0.09210955608663024
$ws.Range("P14").Value = 0.09210955608663024
$ws.Range("Q14").Value = 1.987436005186556
$ws.Range("R14").Value = 17.886924046679
$ws.Range("S14").Value = 0.0001351447931754531
$ws.Range("T14").Value = 0.0001351447931754531
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Apoe"
$ws.Range("C15").Value = "Vldlr"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 6.448603666666666
$ws.Range("H15").Value = 19.345811
$ws.Range("I15").Value = 0.00146721794043115
$ws.Range("J15").Value = 0.00146721794043115
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.207039333333333
$ws.Range("N15").Value = 6.621118
$ws.Range("O15").Value = 0.65961009678592
$ws.Range("P15").Value = 0.6596100967859201
$ws.Range("Q15").Value = 14.23232193741089
$ws.Range("R15").Value = 128.090897436698
$ws.Range("S15").Value = 0.000967791767693829
$ws.Range("T15").Value = 0.0009677917676938293
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Apoe"
$ws.Range("C16").Value = "Vldlr"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 6.448603666666666
$ws.Range("H16").Value = 19.345811
$ws.Range("I16").Value = 0.00146721794043115
$ws.Range("J16").Value = 0.00146721794043115
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.83074
$ws.Range("N16").Value = 2.49222
$ws.Range("O16").Value = 0.2482803471274497
$ws.Range("P16").Value = 0.2482803471274497
$ws.Range("Q16").Value = 5.357113010046667
$ws.Range("R16").Value = 48.21401709041999
$ws.Range("S16").Value = 0.0003642813795618678
$ws.Range("T16").Value = 0.0003642813795618678
